$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (stays visible) ---
$ws.Range("M10").Value = 14.48
$ws.Range("N10").Value = 8.69
$ws.Range("P10").Value = 5
$ws.Range("Q10").Value = 3
$ws.Range("U10").Value = 3

# --- Row 11 (becomes hidden) ---
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("P11").Value = 7
$ws.Range("Q11").Value = 0
$ws.Range("U11").Value = 0
$ws.Rows.Item(11).Hidden = $true

# --- Row 17 (becomes hidden) ---
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("P17").Value = 60
$ws.Range("Q17").Value = 0
$ws.Range("U17").Value = 0
$ws.Rows.Item(17).Hidden = $true

# --- Row 23 (becomes hidden) ---
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("P23").Value = 70
$ws.Range("Q23").Value = 0
$ws.Range("U23").Value = 0
$ws.Rows.Item(23).Hidden = $true

# --- Row 30 (becomes hidden) ---
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("P30").Value = 8
$ws.Range("Q30").Value = 0
$ws.Range("U30").Value = 0
$ws.Rows.Item(30).Hidden = $true

# --- Rows 33 & 34 swap their content (row 33 becomes hidden, row 34 stays hidden) ---
$ws.Range("C33").Value = "M105A20  "
$ws.Range("F33").Value = 14
$ws.Range("G33").Value = 1.72
$ws.Range("H33").Value = 0.6899999999999999
$ws.Range("J33").Value = "REDUCIR 19%"
$ws.Range("K33").Value = 5
$ws.Range("L33").Value = 4
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("P33").Value = 24
$ws.Range("Q33").Value = 0
$ws.Range("S33").Value = 1
$ws.Range("T33").Value = 1
$ws.Range("U33").Value = 0
$ws.Rows.Item(33).Hidden = $true

$ws.Range("C34").Value = "M13A25   "
$ws.Range("F34").Value = 6
$ws.Range("G34").Value = 2.75
$ws.Range("H34").Value = 1.1
$ws.Range("J34").Value = "REDUCIR 9%"
$ws.Range("K34").Value = 2
$ws.Range("L34").Value = 0
$ws.Range("P34").Value = 14
$ws.Range("S34").Value = 0
$ws.Range("T34").Value = 0

# --- Row 37 (becomes hidden) ---
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("P37").Value = 16
$ws.Range("Q37").Value = 0
$ws.Range("U37").Value = 0
$ws.Rows.Item(37).Hidden = $true

# --- Row 43 (becomes hidden) ---
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("P43").Value = 10
$ws.Range("Q43").Value = 0
$ws.Range("U43").Value = 0
$ws.Rows.Item(43).Hidden = $true

# --- Row 44 (becomes hidden) ---
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("P44").Value = 34
$ws.Range("Q44").Value = 0
$ws.Range("U44").Value = 0
$ws.Rows.Item(44).Hidden = $true

# --- Row 46 (becomes hidden) ---
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("P46").Value = 44
$ws.Range("Q46").Value = 0
$ws.Range("U46").Value = 0
$ws.Rows.Item(46).Hidden = $true

# --- Row 54 (becomes hidden) ---
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("P54").Value = 38
$ws.Range("Q54").Value = 0
$ws.Range("U54").Value = 0
$ws.Rows.Item(54).Hidden = $true

# --- Row 68 (stays visible) ---
$ws.Range("L68").Value = 0

# --- Summary rows ---
$ws.Range("C72").Value = 132

# C74 holds a text value that looks numeric/currency ("679.86€"); a direct
# .Value assignment gets auto-parsed by Excel into a number with a new
# currency number format, which would change the cell's style index.
# Instead, write it as a formula that evaluates to the literal text, then
# convert that formula to a static value via copy / paste-values so the
# cell ends up as a plain (shared) string with the original style intact.
$ws.Range("C74").Formula = '="679.86€"'
$ws.Range("C74").Copy() | Out-Null
$ws.Range("C74").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
